$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7
$ws.Range("B7").Value = -0.3456593728802653
$ws.Range("C7").Value = 1.223567895597836
$ws.Range("D7").Value = 3.031172027098284
$ws.Range("E7").Value = 1.741026141991637
$ws.Range("F7").Value = 1.729273255816612
$ws.Range("G7").Value = 38

# Row 8
$ws.Range("B8").Value = -0.3394225716372868
$ws.Range("C8").Value = 1.162219116590446
$ws.Range("D8").Value = 2.462144046345596
$ws.Range("E8").Value = 1.569122062283746
$ws.Range("F8").Value = 1.553103036180247
$ws.Range("G8").Value = 37

# Row 9
$ws.Range("B9").Value = -0.01562264765185476
$ws.Range("C9").Value = 0.5074306157638555
$ws.Range("D9").Value = 0.4207352464598088
$ws.Range("E9").Value = 0.6486410767595657
$ws.Range("F9").Value = 0.6652986502633749
$ws.Range("G9").Value = 20

# Row 10
$ws.Range("B10").Value = 0.1076720123623735
$ws.Range("C10").Value = 0.5887385248463797
$ws.Range("D10").Value = 0.5656078948099947
$ws.Range("E10").Value = 0.7520690758234876
$ws.Range("F10").Value = 0.774714475539742
$ws.Range("G10").Value = 13

# Row 11
$ws.Range("B11").Value = 0.1519798832969088
$ws.Range("C11").Value = 0.4542257695548669
$ws.Range("D11").Value = 0.2920572055036713
$ws.Range("E11").Value = 0.5404231726190794
$ws.Range("F11").Value = 0.5798268282176253
$ws.Range("G11").Value = 5
